$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'54.479.34"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "'2.290.36"
$ws.Range("E3").Value = "  +0.15%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'504.06"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'130.38"
$ws.Range("E6").Value = "  +2.35%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +0.58%  "
$ws.Range("E9").Value = "  +1.91%  "
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("E11").Value = "  +4.98%  "
$ws.Range("D12").Value = "'4.75"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("D13").Value = "'2.698.56"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D15").Value = "'54.452.47"
$ws.Range("E15").Value = "  +0.41%  "
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("E17").Value = "  +0.17%  "
$ws.Range("D18").Value = "'10.26"
$ws.Range("E18").Value = "  +3.20%  "
$ws.Range("D20").Value = "'304.78"
$ws.Range("E20").Value = "  +0.88%  "
$ws.Range("D21").Value = "'6.39"
$ws.Range("E21").Value = "  -0.67%  "
$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.12%  "
$ws.Range("D23").Value = "'61.99"
$ws.Range("E23").Value = "  -2.78%  "
$ws.Range("D24").Value = "'0.996"
$ws.Range("E24").Value = "  -0.69%  "
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("E26").Value = "  +3.72%  "
$ws.Range("D27").Value = "'172.03"
$ws.Range("E27").Value = "  +3.26%  "
$ws.Range("E28").Value = "  +2.18%  "
$ws.Range("D29").Value = "'0.0₃0694"
$ws.Range("E29").Value = "  +1.59%  "
$ws.Range("D30").Value = "'5.98"
$ws.Range("E30").Value = "  +1.57%  "
$ws.Range("D31").Value = "'1.10"
$ws.Range("E31").Value = "  +1.98%  "
$ws.Range("E32").Value = "  -0.05%  "
$ws.Range("D33").Value = "'17.89"
$ws.Range("E33").Value = "  +1.84%  "
$ws.Range("D34").Value = "'0.963"
$ws.Range("E34").Value = "  +10.49%  "
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").Value = "'1.20"
$ws.Range("E36").Value = "  +1.28%  "
$ws.Range("D37").Value = "'3.74"
$ws.Range("E37").Value = "  +3.39%  "
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("E39").Value = "  +1.48%  "
$ws.Range("E40").Value = "  +1.72%  "
$ws.Range("D41").Value = "'4.89"
$ws.Range("E41").Value = "  +1.92%  "
$ws.Range("D42").Value = "'126.07"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("E43").Value = "  +3.84%  "
$ws.Range("E44").Value = "  +1.02%  "
$ws.Range("D45").Value = "'0.550"
$ws.Range("E45").Value = "  +1.22%  "
$ws.Range("D46").Value = "'243.19"
$ws.Range("E46").Value = "  +1.83%  "
$ws.Range("E47").Value = "  +0.49%  "
$ws.Range("D48").Value = "'0.0206"
$ws.Range("E48").Value = "  +1.85%  "
$ws.Range("E49").Value = "  +0.81%  "
$ws.Range("D50").Value = "'16.48"
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("E51").Value = "  +1.84%  "
